# "Commercial in One Flow" — refresh the Application No / Consumer Number
# test-data values on the "Commercial" sheet for the next automation run.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Commercial")

# AN1/AO1 hold the "Application No" / "Consumer Number" headers (unchanged);
# AN2/AO2 hold the actual test values that need bumping to the new numbers.
$ws.Range("AN2").Value = "JP30000195"
$ws.Range("AO2").Value = "JP30000097"

# Reflect the user having the two trailing columns selected/focused afterwards.
$ws.Range("AN1:AO1048576").Select()
